$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the LED "Value" description with a clarifying comment explaining
# the correct orientation of the SK6812 3535 LED chip on the PCB.
$ws.Range("C4").Value = "SK6812 3535 (chip inside must be facing side with no leads - up on the PCB)"

# Enable word-wrap for the updated cell and grow the row so the longer text is readable.
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 25.2

# Move the active selection (this also clears the old frozen/scrolled topLeftCell).
$ws.Range("G13").Select()
